$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Fri Oct 06 11:30:14 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:30:27 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:30:41 EDT 2023"
$ws.Range("B5").Value = "Fri Oct 06 11:30:54 EDT 2023"
